# Update the "dSF" (column F) values for the affected rows to reflect the
# repulled data / pushed data / mean calculation referenced in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -1
    7  = -2
    9  = 0
    10 = -2
    11 = -4
    12 = -4
    15 = -3
    16 = 0
    20 = -5
    22 = -4
    30 = 3
    31 = -2
    32 = 1
    33 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
